# Weekly refresh of the Fruta/Hortaliza "Tuna" dataset: the rows of
# observations (columns D, M, N, O, P, Q, R, S, T) get re-shuffled across
# dates (row 8 is untouched). We snapshot every source row first so the
# row-to-row copy is safe regardless of write order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# target row -> source row (values for target row are taken from source row)
$mapping = @{
    2  = 3
    3  = 13
    4  = 18
    5  = 4
    6  = 17
    7  = 20
    9  = 6
    10 = 14
    11 = 22
    12 = 7
    13 = 5
    14 = 10
    15 = 2
    16 = 15
    17 = 11
    18 = 21
    19 = 9
    20 = 19
    21 = 12
    22 = 16
}

# Columns touched by the refresh.
$cols = @(4, 13, 14, 15, 16, 17, 18, 19, 20)

# Snapshot the current ("before") value of every touched cell so that
# overwriting earlier rows doesn't clobber data still needed for later rows.
$snapshot = @{}
foreach ($r in $mapping.Keys) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Write each target row's values from its mapped source row's snapshot.
foreach ($r in $mapping.Keys) {
    $src = $mapping[$r]
    $srcVals = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value = $srcVals[$c]
    }
}
